$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("register")

# Update first/last-name style fields (column A)
$ws.Range("A2").Value = "dddhivya"
$ws.Range("A3").Value = "adanjum"
$ws.Range("A4").Value = "sfshivender"
$ws.Range("A5").Value = "nfnaveen"

# Update email fields (column C)
$ws.Range("C2").Value = "dhivya111abc@gmail.com"
$ws.Range("C3").Value = "anjum111abc@gmail.com"
$ws.Range("C4").Value = "shiv112abc@gmail.com"
$ws.Range("C5").Value = "naveen11a2bc@gmail.com"
